# TMX IAMC Indigenous Community Profiles workbook update
# - adds a new "Spreads" worksheet with a Number/Sub/Start/Stop lookup table
#   (prep work for "kp click" functionality)
# - widens column O on "BC First Nations" to fit the new data and moves the
#   active cell / active tab there
# - recalculation of the "image sources" lookups happens automatically

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Add the new "Spreads" sheet as the last tab
# ---------------------------------------------------------------------
$sheetCount = $wb.Worksheets.Count
$ws = $wb.Worksheets.Add($null, $wb.Worksheets.Item($sheetCount))
$ws.Name = "Spreads"

# Header row
$ws.Range("A1").Value = "Number"
$ws.Range("B1").Value = "Sub"
$ws.Range("C1").Value = "Start"
$ws.Range("D1").Value = "Stop"

# Data rows: Number, Sub (blank for some rows), Start, Stop
$data = @(
  @(1, "", 0, 49.114),
  @(2, "A", 49.114, 147.767),
  @(2, "B", 147.767, 246.043),
  @(2, "C", 246.043, 337.966),
  @(3, "A", 489.16, 502.443),
  @(3, "B", 502.443, 525.579),
  @(3, "C", 525.579, 610.675),
  @(4, "A", 610.675, 690.485),
  @(4, "B", 690.485, 764.569),
  @(5, "A", 806.471, 990.273),
  @(5, "B", 990.273, 1075.731),
  @(6, "", 1075.122, 1144.48),
  @(7, "", 1144.908, 1180.149)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    if ($row[1] -ne "") {
        $ws.Cells.Item($r, 2).Value = $row[1]
    }
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

$ws.Range("A15").Select()

# ---------------------------------------------------------------------
# 2. Widen column O on "BC First Nations" and move the active cell there
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("BC First Nations")
$ws1.Select()
$ws1.Columns.Item(15).ColumnWidth = 40
$ws1.Range("O12").Select()

# ---------------------------------------------------------------------
# 3. Leave the new "Spreads" sheet as the active tab
# ---------------------------------------------------------------------
$ws.Select()
